$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new sale/item line -------------------------------------------
# A7: index number 0 -> 1
$ws.Cells.Item(7, 1).Value = 1

# C7:G7 and N7:O7 share one visual style; switch it to Text (@) BEFORE
# writing the strings so they are stored as literal text, not numbers.
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "BI ALCOFAN 150 MG 30 TABS."   # C7 - item name
$ws.Cells.Item(7, 14).Value = "81.00"                        # N7 - price

# H7:K7 share another style; switch to Text as well.
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Cells.Item(7, 8).Value = "2:2"                            # H7 - current balance

# Q7 keeps its own style; switch to Text.
$ws.Range("Q7").NumberFormat = "@"
$ws.Cells.Item(7, 17).Value = "0:1"                           # Q7 - transaction count

# L7 and P7 must stay on their ORIGINAL number formats (custom #,##0.##
# and 0.00 respectively) while still holding literal text. Flip the cell
# to Text only long enough to assign the string, then restore the
# original format so the stored cellXfs entry is left untouched.
$fmtL7 = $ws.Cells.Item(7, 12).NumberFormat
$ws.Cells.Item(7, 12).NumberFormat = "@"
$ws.Cells.Item(7, 12).Value = "1"                             # L7 - order limit
$ws.Cells.Item(7, 12).NumberFormat = $fmtL7

$fmtP7 = $ws.Cells.Item(7, 16).NumberFormat
$ws.Cells.Item(7, 16).NumberFormat = "@"
$ws.Cells.Item(7, 16).Value = "26.7300"                       # P7 - sale price
$ws.Cells.Item(7, 16).NumberFormat = $fmtP7

# --- Row 8: sale-price subtotal -------------------------------------------
$ws.Cells.Item(8, 16).Value = 26.73                           # P8

# --- Row 9: refreshed report timestamp ------------------------------------
$ws.Cells.Item(9, 1).Value = "Wednesday, 24 September, 2025 10:13 AM"
